$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3955.5557
$ws.Range("I64").Value = 2900
$ws.Range("J64").Value = 4257.143
$ws.Range("K64").Value = 2900
$ws.Range("L64").Value = 4257.143
$ws.Range("M64").Value = -2652
$ws.Range("N64").Value = -4753.143
$ws.Range("H67").Value = 3955.5557
$ws.Range("I67").Value = 2900
$ws.Range("J67").Value = 4257.143
$ws.Range("K67").Value = 2900
$ws.Range("L67").Value = 4257.143
$ws.Range("M67").Value = -2042
$ws.Range("N67").Value = -5973.143
$ws.Range("H106").Value = 2862.9707
$ws.Range("I106").Value = 1794.7858
$ws.Range("K106").Value = 1794.7858
$ws.Range("M106").Value = -1163.7858
$ws.Range("H137").Value = 107891.84
$ws.Range("I137").Value = 135409.67
$ws.Range("K137").Value = 406229.01
$ws.Range("M137").Value = -403679.01
$ws.Range("H138").Value = 3442.1467
$ws.Range("J138").Value = 3379.1875
$ws.Range("L138").Value = 10137.5625
$ws.Range("N138").Value = -20417.5625

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14172.743
$ws.Range("I32").Value = 9903.299999999999
$ws.Range("K32").Value = 9903.299999999999
$ws.Range("M32").Value = -9616.299999999999
$ws.Range("H61").Value = 18495302
$ws.Range("I61").Value = 29281896
$ws.Range("K61").Value = 29281896
$ws.Range("M61").Value = -29281684
$ws.Range("H74").Value = 31251864
$ws.Range("I74").Value = 52632376
$ws.Range("K74").Value = 52632376
$ws.Range("M74").Value = -52631502
$ws.Range("H77").Value = 31251864
$ws.Range("I77").Value = 52632376
$ws.Range("K77").Value = 263161880
$ws.Range("M77").Value = -263157512
$ws.Range("H132").Value = 10883344
$ws.Range("I132").Value = 14288236
$ws.Range("K132").Value = 42864708
$ws.Range("M132").Value = -42862178
$ws.Range("H136").Value = 18495302
$ws.Range("I136").Value = 29281896
$ws.Range("K136").Value = 87845688
$ws.Range("M136").Value = -87843138
$ws.Range("H139").Value = 42690.8
$ws.Range("J139").Value = 42690.8
$ws.Range("L139").Value = 42690.8
$ws.Range("N139").Value = -52970.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H105").Value = 5858017
$ws.Range("I105").Value = 11906344
$ws.Range("J105").Value = 2176426.5
$ws.Range("K105").Value = 11906344
$ws.Range("L105").Value = 2176426.5
$ws.Range("M105").Value = -11904597
$ws.Range("N105").Value = -2179920.5
$ws.Range("H134").Value = 3212.8928
$ws.Range("I134").Value = 2996.3914
$ws.Range("J134").Value = 4208.8
$ws.Range("K134").Value = 8989.174199999999
$ws.Range("L134").Value = 12626.4
$ws.Range("M134").Value = -6454.174199999999
$ws.Range("N134").Value = -17696.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5711.778
$ws.Range("I62").Value = 4733.3335
$ws.Range("K62").Value = 4733.3335
$ws.Range("M62").Value = -4109.3335
$ws.Range("H65").Value = 5711.778
$ws.Range("I65").Value = 4733.3335
$ws.Range("K65").Value = 23666.6675
$ws.Range("M65").Value = -20546.6675

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1335.4286
$ws.Range("I5").Value = 983.5806
$ws.Range("J5").Value = 2327
$ws.Range("K5").Value = 2950.7418
$ws.Range("L5").Value = 6981
$ws.Range("M5").Value = -2838.7418
$ws.Range("N5").Value = -7205
$ws.Range("H131").Value = 714.74
$ws.Range("I131").Value = 415.45456
$ws.Range("J131").Value = 751.73035
$ws.Range("K131").Value = 1246.36368
$ws.Range("L131").Value = 2255.19105
$ws.Range("M131").Value = 3793.63632
$ws.Range("N131").Value = -12335.19105
$ws.Range("H134").Value = 4244.952
$ws.Range("I134").Value = 3213
$ws.Range("J134").Value = 5183.091
$ws.Range("K134").Value = 9639
$ws.Range("L134").Value = 15549.273
$ws.Range("M134").Value = -4569
$ws.Range("N134").Value = -25689.273
$ws.Range("H135").Value = 1335.4286
$ws.Range("I135").Value = 983.5806
$ws.Range("J135").Value = 2327
$ws.Range("K135").Value = 8852.225399999999
$ws.Range("L135").Value = 20943
$ws.Range("M135").Value = -6317.225399999999
$ws.Range("N135").Value = -26013
$ws.Range("H136").Value = 3553.3333
$ws.Range("J136").Value = 4396
$ws.Range("L136").Value = 13188
$ws.Range("N136").Value = -23388
$ws.Range("H141").Value = 2015
$ws.Range("I141").Value = 2015
$ws.Range("K141").Value = 6045
$ws.Range("M141").Value = -865

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3994.7222
$ws.Range("I80").Value = 3724.875
$ws.Range("J80").Value = 4210.6
$ws.Range("K80").Value = 3724.875
$ws.Range("L80").Value = 4210.6
$ws.Range("M80").Value = -2726.875
$ws.Range("N80").Value = -6206.6
$ws.Range("H83").Value = 3994.7222
$ws.Range("I83").Value = 3724.875
$ws.Range("J83").Value = 4210.6
$ws.Range("K83").Value = 18624.375
$ws.Range("L83").Value = 21053
$ws.Range("M83").Value = -13632.375
$ws.Range("N83").Value = -31037
$ws.Range("H97").Value = 2701.4285
$ws.Range("I97").Value = 2533
$ws.Range("K97").Value = 2533
$ws.Range("M97").Value = -2037
$ws.Range("H132").Value = 3190856
$ws.Range("I132").Value = 3971900.8
$ws.Range("K132").Value = 11915702.4
$ws.Range("M132").Value = -11913172.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4142.4546
$ws.Range("I7").Value = 3429.111
$ws.Range("J7").Value = 7352.5
$ws.Range("K7").Value = 3429.111
$ws.Range("L7").Value = 7352.5
$ws.Range("M7").Value = -3317.111
$ws.Range("N7").Value = -7576.5
$ws.Range("H93").Value = 1556.3125
$ws.Range("I93").Value = 1556.3125
$ws.Range("K93").Value = 1556.3125
$ws.Range("M93").Value = -308.3125
$ws.Range("H126").Value = 4142.4546
$ws.Range("I126").Value = 3429.111
$ws.Range("J126").Value = 7352.5
$ws.Range("K126").Value = 10287.333
$ws.Range("L126").Value = 22057.5
$ws.Range("M126").Value = -7817.332999999999
$ws.Range("N126").Value = -26997.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H47").Value = 9434.5
$ws.Range("J47").Value = 9434.5
$ws.Range("L47").Value = 9434.5
$ws.Range("N47").Value = -10578.5
$ws.Range("H122").Value = 2099.875
$ws.Range("I122").Value = 1966.6666
$ws.Range("J122").Value = 2499.5
$ws.Range("K122").Value = 5899.9998
$ws.Range("L122").Value = 7498.5
$ws.Range("M122").Value = -3449.9998
$ws.Range("N122").Value = -12398.5
$ws.Range("H126").Value = 2559.3125
$ws.Range("I126").Value = 2144.4
$ws.Range("J126").Value = 3250.8333
$ws.Range("K126").Value = 6433.200000000001
$ws.Range("L126").Value = 9752.499899999999
$ws.Range("M126").Value = -3963.200000000001
$ws.Range("N126").Value = -14692.4999
$ws.Range("H132").Value = 15153169
$ws.Range("I132").Value = 20834606
$ws.Range("J132").Value = 2672
$ws.Range("K132").Value = 62503818
$ws.Range("L132").Value = 8016
$ws.Range("M132").Value = -62501288
$ws.Range("N132").Value = -13076
$ws.Range("H136").Value = 26318844
$ws.Range("I136").Value = 45456296
$ws.Range("J136").Value = 4844.0625
$ws.Range("K136").Value = 136368888
$ws.Range("L136").Value = 14532.1875
$ws.Range("M136").Value = -136366338
